# Exercise with FEM in tools folder - update element table values and
# tidy up the sheet layout/view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Cell value updates (Capacity / Conductivity columns) -----------------
# Capacity (column C) for the Top/Middle/Bottom rows
$ws.Range("C4").Value = 100
$ws.Range("C5").Value = 200
$ws.Range("C6").Value = 300

# Conductivity (column D) for rows 4-10
$ws.Range("D4").Value = 1
$ws.Range("D5").Value = 2
$ws.Range("D6").Value = 3
$ws.Range("D7").Value = 4
$ws.Range("D8").Value = 5
$ws.Range("D9").Value = 6
$ws.Range("D10").Value = 7

# --- Column widths ----------------------------------------------------
# Target stored widths are 21.140625 / 17.140625 / 16 "characters".
# ColumnWidth in this host snaps to the nearest 1/6 of a character, so use
# the values that land on (or closest to) those stored widths.
$ws.Columns.Item(2).ColumnWidth = 20.334
$ws.Columns.Item(3).ColumnWidth = 16.334
$ws.Columns.Item(4).ColumnWidth = 15.167

# --- Selection ----------------------------------------------------------
$ws.Range("K10").Select()

# --- Window geometry ------------------------------------------------------
$win = $excel.ActiveWindow
$win.Left = 1455
$win.Top = 765
$win.Width = 16845
$win.Height = 11295
